$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted before the existing row 181
# (Fecha 2022-11-08 / serial 44873), pushing every subsequent row (old 181..290)
# down by one (new 182..291). Insert a fresh row at 181 first so the rows
# below shift down automatically, then populate the new row's values.
$ws.Rows.Item(181).Insert()

$ws.Range("A181").Value = 7
$ws.Range("B181").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C181").Value = "Ñuble"
$ws.Range("D181").Value = 44873
$ws.Range("E181").Value = 16
$ws.Range("F181").Value = 100112003
$ws.Range("G181").Value = "Ajo"
$ws.Range("H181").Value = "Chino"
$ws.Range("I181").Value = "Primera"
$ws.Range("J181").Value = 60
$ws.Range("K181").Value = 17000
$ws.Range("L181").Value = 18000
$ws.Range("M181").Value = 17500
$ws.Range("N181").Value = "$/malla 10 kilos"
$ws.Range("O181").Value = "China"
$ws.Range("P181").Value = 1750
$ws.Range("Q181").Value = 10
$ws.Range("R181").Value = "Hortaliza"
